$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D15").Value = 7
$ws.Range("E15").Value = 7
$ws.Range("F15").Value = 10

$ws.Range("E16").Value = 9
$ws.Range("F16").Value = 9

$ws.Range("E17").Value = 10
$ws.Range("F17").Value = 10

$ws.Range("E18").Value = 10
$ws.Range("F18").Value = 10

$ws.Range("E19").Value = 10
$ws.Range("F19").Value = 9

$ws.Range("D20").Value = 9
$ws.Range("E20").Value = 9
$ws.Range("F20").Value = 8

$ws.Range("F20").Select()
